$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(2, 2).Range.Text = "349"
$t.Cell(2, 3).Range.Text = "157 (45.0)"

$t.Cell(3, 2).Range.Text = "106"
$t.Cell(3, 3).Range.Text = "47 (44.3)"

$t.Cell(4, 2).Range.Text = "109"
$t.Cell(4, 3).Range.Text = "51 (46.8)"

$t.Cell(5, 2).Range.Text = "107"
$t.Cell(5, 3).Range.Text = "43 (40.2)"

$t.Cell(6, 2).Range.Text = "27"
$t.Cell(6, 3).Range.Text = "16 (59.3)"
